$d = $word.ActiveDocument

# Locate the "Author" styled paragraph that holds "Edison Achalma" (the
# byline under the title) and append a new "Author" paragraph right
# after it with the institutional affiliation.
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Author" -and $p.Range.Text.Trim() -eq "Edison Achalma") {
        $p.Range.InsertAfter([char]13 + "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga")
        break
    }
}
